$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics for Fndc5-Itgav ligand-receptor pairs.
# Values below are the new literal numbers written by the upstream Python pipeline;
# Sending cluster columns (E:J) and Target cluster columns (K:P) are recomputed per
# the refreshed TPM expression values, with the Edge columns (Q:T) following from them.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.685737
$ws.Range("H2").Value = 2.057211
$ws.Range("I2").Value = 0.05519567570004053
$ws.Range("J2").Value = 0.05519567570004053
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 6.462284227503001
$ws.Range("R2").Value = 58.160558047527
$ws.Range("S2").Value = 0.00367290587873599
$ws.Range("T2").Value = 0.00367290587873599

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.685737
$ws.Range("H3").Value = 2.057211
$ws.Range("I3").Value = 0.05519567570004053
$ws.Range("J3").Value = 0.05519567570004053
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 34.693193516826
$ws.Range("R3").Value = 312.238741651434
$ws.Range("S3").Value = 0.01971823428591471
$ws.Range("T3").Value = 0.01971823428591471

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.685737
$ws.Range("H4").Value = 2.057211
$ws.Range("I4").Value = 0.05519567570004053
$ws.Range("J4").Value = 0.05519567570004053
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 18.405704068752
$ws.Range("R4").Value = 165.651336618768
$ws.Range("S4").Value = 0.01046107170413262
$ws.Range("T4").Value = 0.01046107170413262

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.685737
$ws.Range("H5").Value = 2.057211
$ws.Range("I5").Value = 0.05519567570004053
$ws.Range("J5").Value = 0.05519567570004053
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 37.552699206243
$ws.Range("R5").Value = 337.974292856187
$ws.Range("S5").Value = 0.02134346383125721
$ws.Range("T5").Value = 0.02134346383125721

# Row 6
$ws.Range("I6").Value = 0.2871009238089374
$ws.Range("J6").Value = 0.2871009238089374
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 33.61364360706023
$ws.Range("R6").Value = 302.522792463542
$ws.Range("S6").Value = 0.01910466096255445
$ws.Range("T6").Value = 0.01910466096255445

# Row 7
$ws.Range("I7").Value = 0.2871009238089374
$ws.Range("J7").Value = 0.2871009238089374
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.1025646159335453
$ws.Range("T7").Value = 0.1025646159335453

# Row 8
$ws.Range("I8").Value = 0.2871009238089374
$ws.Range("J8").Value = 0.2871009238089374
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 95.73747534516978
$ws.Range("R8").Value = 861.6372781065281
$ws.Range("S8").Value = 0.05441338134186123
$ws.Range("T8").Value = 0.05441338134186125

# Row 9
$ws.Range("I9").Value = 0.2871009238089374
$ws.Range("J9").Value = 0.2871009238089374
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 195.3307844662113
$ws.Range("R9").Value = 1757.977060195902
$ws.Range("S9").Value = 0.1110182655709764
$ws.Range("T9").Value = 0.1110182655709764

# Row 10
$ws.Range("G10").Value = 7.915626666666667
$ws.Range("H10").Value = 23.74688
$ws.Range("I10").Value = 0.6371369234209706
$ws.Range("J10").Value = 0.6371369234209707
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 74.59569683246224
$ws.Range("R10").Value = 671.3612714921601
$ws.Range("S10").Value = 0.04239723351354727
$ws.Range("T10").Value = 0.04239723351354728

# Row 11
$ws.Range("G11").Value = 7.915626666666667
$ws.Range("H11").Value = 23.74688
$ws.Range("I11").Value = 0.6371369234209706
$ws.Range("J11").Value = 0.6371369234209707
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 400.4718540105245
$ws.Range("R11").Value = 3604.24668609472
$ws.Range("S11").Value = 0.2276123078281724
$ws.Range("T11").Value = 0.2276123078281724

# Row 12
$ws.Range("G12").Value = 7.915626666666667
$ws.Range("H12").Value = 23.74688
$ws.Range("I12").Value = 0.6371369234209706
$ws.Range("J12").Value = 0.6371369234209707
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 212.4614567179378
$ws.Range("R12").Value = 1912.15311046144
$ws.Range("S12").Value = 0.1207546597939797
$ws.Range("T12").Value = 0.1207546597939797

# Row 13
$ws.Range("G13").Value = 7.915626666666667
$ws.Range("H13").Value = 23.74688
$ws.Range("I13").Value = 0.6371369234209706
$ws.Range("J13").Value = 0.6371369234209707
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 433.4798140427733
$ws.Range("R13").Value = 3901.31832638496
$ws.Range("S13").Value = 0.2463727222852713
$ws.Range("T13").Value = 0.2463727222852713

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2555126666666667
$ws.Range("H14").Value = 0.7665379999999999
$ws.Range("I14").Value = 0.02056647707005147
$ws.Range("J14").Value = 0.02056647707005147
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 2.407913639962889
$ws.Range("R14").Value = 21.671222759666
$ws.Range("S14").Value = 0.001368562547290739
$ws.Range("T14").Value = 0.001368562547290739

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2555126666666667
$ws.Range("H15").Value = 0.7665379999999999
$ws.Range("I15").Value = 0.02056647707005147
$ws.Range("J15").Value = 0.02056647707005147
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 12.92704111148578
$ws.Range("R15").Value = 116.343370003372
$ws.Range("S15").Value = 0.007347217117279896
$ws.Range("T15").Value = 0.007347217117279896

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2555126666666667
$ws.Range("H16").Value = 0.7665379999999999
$ws.Range("I16").Value = 0.02056647707005147
$ws.Range("J16").Value = 0.02056647707005147
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 6.858154844327111
$ws.Range("R16").Value = 61.723393598944
$ws.Range("S16").Value = 0.003897903025962046
$ws.Range("T16").Value = 0.003897903025962047

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2555126666666667
$ws.Range("H17").Value = 0.7665379999999999
$ws.Range("I17").Value = 0.02056647707005147
$ws.Range("J17").Value = 0.02056647707005147
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 13.99252237332733
$ws.Range("R17").Value = 125.932701359946
$ws.Range("S17").Value = 0.007952794379518794
$ws.Range("T17").Value = 0.007952794379518794
